# Transportation sector calibration and redesignating passenger ships as taxis
$wb = $excel.ActiveWorkbook

$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

# Row 2 - LDVs: 0.060966558335291694 -> 0.051
$wsPsgr.Range("B2:H2").Value = 0.051

# Row 3 - HDVs: 0.1 -> 0.096
$wsPsgr.Range("B3:H3").Value = 0.096

# Row 4 - aircraft: 0.05 -> 0.045
$wsPsgr.Range("B4:H4").Value = 0.045

# Row 6 - ships: 0.030303 -> 0.3 (redesignating passenger ships as taxis)
$wsPsgr.Range("B6:H6").Value = 0.3

# Cosmetic / view updates to match target workbook
$wsPsgr.Range("B6").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("G3").Select()

$wb.Save()
